# Swap the order of "<email>, System" -> "System, <email>" in column G
# (Recorded By) for every applicable row, except entries for
# backup@backdoor.com which keep their original order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

$suffix = ", System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    if ($val.EndsWith($suffix)) {
        $email = $val.Substring(0, $val.Length - $suffix.Length)
        if ($email -ne "backup@backdoor.com") {
            $cell.Value = "System, " + $email
        }
    }
}
